$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 2840.2
$ws.Range("I2").Value = 5300.5
$ws.Range("K2").Value = 5300.5
$ws.Range("M2").Value = -5187.5
$ws.Range("H13").Value = 3329.9
$ws.Range("I13").Value = 2175
$ws.Range("J13").Value = 4099.8335
$ws.Range("K13").Value = 2175
$ws.Range("L13").Value = 4099.8335
$ws.Range("M13").Value = -2006
$ws.Range("N13").Value = -4437.8335
$ws.Range("H16").Value = 2750
$ws.Range("I16").Value = 2500
$ws.Range("K16").Value = 2500
$ws.Range("M16").Value = -2270
$ws.Range("H18").Value = 3341.5
$ws.Range("J18").Value = 3199.5
$ws.Range("L18").Value = 3199.5
$ws.Range("N18").Value = -3767.5
$ws.Range("H43").Value = 12953.066
$ws.Range("I43").Value = 12916.333
$ws.Range("J43").Value = 12977.556
$ws.Range("K43").Value = 12916.333
$ws.Range("L43").Value = 12977.556
$ws.Range("M43").Value = -12847.333
$ws.Range("N43").Value = -13115.556
$ws.Range("H48").Value = 1499
$ws.Range("I48").Value = 1499
$ws.Range("K48").Value = 4497
$ws.Range("M48").Value = -4205
$ws.Range("H56").Value = 1499
$ws.Range("I56").Value = 1499
$ws.Range("K56").Value = 4497
$ws.Range("M56").Value = -3963
$ws.Range("H136").Value = 84249.75
$ws.Range("J136").Value = 84249.75
$ws.Range("L136").Value = 84249.75
$ws.Range("N136").Value = -94449.75
$ws.Range("H137").Value = 5384.8
$ws.Range("J137").Value = 20002.5
$ws.Range("L137").Value = 60007.5
$ws.Range("N137").Value = -65107.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8353254
$ws.Range("I32").Value = 10225005
$ws.Range("J32").Value = 15457.091
$ws.Range("K32").Value = 10225005
$ws.Range("L32").Value = 15457.091
$ws.Range("M32").Value = -10224718
$ws.Range("N32").Value = -16031.091
$ws.Range("H45").Value = 25002088
$ws.Range("I45").Value = 35716124
$ws.Range("K45").Value = 35716124
$ws.Range("M45").Value = -35715747
$ws.Range("H74").Value = 22739628
$ws.Range("J74").Value = 18416.572
$ws.Range("L74").Value = 18416.572
$ws.Range("N74").Value = -20164.572
$ws.Range("H77").Value = 22739628
$ws.Range("J77").Value = 18416.572
$ws.Range("L77").Value = 92082.86
$ws.Range("N77").Value = -100818.86

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2604.0454
$ws.Range("I99").Value = 1498.4445
$ws.Range("K99").Value = 1498.4445
$ws.Range("M99").Value = -0.4445000000000618
$ws.Range("H105").Value = 2428.5806
$ws.Range("I105").Value = 1580.4
$ws.Range("J105").Value = 2832.476
$ws.Range("K105").Value = 1580.4
$ws.Range("L105").Value = 2832.476
$ws.Range("M105").Value = 166.5999999999999
$ws.Range("N105").Value = -6326.476000000001
$ws.Range("H107").Value = 3662.7368
$ws.Range("I107").Value = 3340.6
$ws.Range("K107").Value = 3340.6
$ws.Range("M107").Value = -1420.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 2500
$ws.Range("J4").Value = 2500
$ws.Range("L4").Value = 2500
$ws.Range("N4").Value = -2724
$ws.Range("H16").Value = 1099.75
$ws.Range("I16").Value = 1066.3334
$ws.Range("K16").Value = 1066.3334
$ws.Range("M16").Value = -779.3334
$ws.Range("H22").Value = 500
$ws.Range("I22").Value = 500
$ws.Range("K22").Value = 500
$ws.Range("M22").Value = -150
$ws.Range("H23").Value = 30000
$ws.Range("J23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("N23").ClearContents()
$ws.Range("H27").Value = 30000
$ws.Range("J27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("N27").ClearContents()
$ws.Range("H31").Value = 1308161.5
$ws.Range("J31").Value = 1798838.9
$ws.Range("L31").Value = 1798838.9
$ws.Range("N31").Value = -1799428.9
$ws.Range("H34").Value = 1308161.5
$ws.Range("J34").Value = 1798838.9
$ws.Range("L34").Value = 1798838.9
$ws.Range("N34").Value = -1799242.9
$ws.Range("H113").Value = 1099.75
$ws.Range("I113").Value = 1066.3334
$ws.Range("K113").Value = 1066.3334
$ws.Range("M113").Value = 1103.6666

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 15551
$ws.Range("J131").Value = 18373.455
$ws.Range("L131").Value = 55120.36500000001
$ws.Range("N131").Value = -65200.36500000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 202.53847
$ws.Range("I2").Value = 54.4
$ws.Range("J2").Value = 295.125
$ws.Range("K2").Value = 54.4
$ws.Range("L2").Value = 295.125
$ws.Range("M2").Value = 58.6
$ws.Range("N2").Value = -521.125
$ws.Range("H107").Value = 0
$ws.Range("I107").Value = 0
$ws.Range("K107").Value = 0
$ws.Range("M107").ClearContents()
$ws.Range("H132").Value = 142858370
$ws.Range("I132").Value = 250001010
$ws.Range("K132").Value = 750003030
$ws.Range("M132").Value = -750000500
$ws.Range("H135").Value = 80000
$ws.Range("J135").Value = 80000
$ws.Range("L135").Value = 80000
$ws.Range("N135").Value = -90140

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H20").Value = 90006
$ws.Range("J20").Value = 90006
$ws.Range("L20").Value = 90006
$ws.Range("N20").Value = -90458
$ws.Range("H22").Value = 1985.8182
$ws.Range("I22").Value = 2126.6365
$ws.Range("J22").Value = 1845
$ws.Range("K22").Value = 2126.6365
$ws.Range("L22").Value = 1845
$ws.Range("M22").Value = -1831.6365
$ws.Range("N22").Value = -2435
$ws.Range("H27").Value = 1985.8182
$ws.Range("I27").Value = 2126.6365
$ws.Range("J27").Value = 1845
$ws.Range("K27").Value = 2126.6365
$ws.Range("L27").Value = 1845
$ws.Range("M27").Value = -2019.6365
$ws.Range("N27").Value = -2059
$ws.Range("H40").Value = 5408.2104
$ws.Range("I40").Value = 5195.3335
$ws.Range("K40").Value = 5195.3335
$ws.Range("M40").Value = -5059.3335

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 21740272
$ws.Range("I107").Value = 27778822
$ws.Range("K107").Value = 83336466
$ws.Range("M107").Value = -83334546
$ws.Range("H132").Value = 504483.3
$ws.Range("I132").Value = 3554.2144
$ws.Range("K132").Value = 10662.6432
$ws.Range("M132").Value = -8132.643199999999
$ws.Range("H136").Value = 6491.9443
$ws.Range("I136").Value = 5877.357
$ws.Range("J136").Value = 8643
$ws.Range("K136").Value = 17632.071
$ws.Range("L136").Value = 25929
$ws.Range("M136").Value = -15082.071
$ws.Range("N136").Value = -31029
